$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "matt_karting"
$ws.Range("B5").Value = "a"
$ws.Range("C5").Value = "a"
$ws.Range("D5").Value = "a"
